$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts existing rows 7-39 down to 8-40,
# matching the diff's row-shift pattern). The new row starts out as a copy
# of the formatting of the surrounding rows (date style on column D, etc.)
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with a new weekly record. Most columns
# (Mercado ID, Mercado, Región, Codreg, Tipo, Producto ID, Producto,
# Categoría ID, Categoría, Variedad, Calidad, Unidad de comercialización,
# Origen, Kg/unidad) repeat the same values used throughout this sheet for
# "Femacal de La Calera" / Papaya "Primera" with the 10kg tray unit.
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 'Femacal de La Calera'
$ws.Range("C7").Value = 'Coquimbo'
$ws.Range("D7").Value = 44473
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 'Fruta'
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = 'Tropicales y subtropicales'
$ws.Range("I7").Value = 100108004
$ws.Range("J7").Value = 'Papaya'
$ws.Range("K7").Value = 'Cultivar IV Región'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 85
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 18000
$ws.Range("P7").Value = 18000
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("R7").Value = 'Provincia del Elquí'
$ws.Range("S7").Value = 1800
$ws.Range("T7").Value = 10
